# Clean up sample order data: clear out a handful of cells that held
# duplicate/placeholder values so the "empty" error-handling paths can be
# exercised, and drop the now-unused "Santas Workshop - Deluxe Edition"
# product row's name value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: clear the duplicate "item" value.
$ws.Range("C2").ClearContents()

# Row 3: clear the duplicate "holiday" value and the product name
# (this was the only remaining use of "Santas Workshop - Deluxe Edition",
# so it drops out of the shared strings table entirely on save).
$ws.Range("B3").ClearContents()
$ws.Range("D3").ClearContents()

# Reset the view: scroll back to the top-left (A1) instead of being
# scrolled over to column I, and move the active selection to D3.
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 1
$aw.ScrollRow = 1
$ws.Range("D3").Select()
